$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cronograma")
$ws.Range("D16").Value = "Criação do banco e tabelas de animais, doações, adm e mensagens"
